$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"
$zhDetail  = "Handback file name: t52hqlad.pst is different with handoff file name: 8548ef9b-ac04-40b2-8971-567febf762a1.2002fd2393f34839bcc3ccc4288c44cd11e08715.zh-cn."
$deDetail  = "Handback file name: t52hqlad.pst is different with handoff file name: 8548ef9b-ac04-40b2-8971-567febf762a1.2002fd2393f34839bcc3ccc4288c44cd11e08715.de-de."

# Overview sheet: row 3 is the 8548ef9b-...md file; update its zh-cn / de-de status cells
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# zh-cn detail sheet: row 3 status + new Error Detail cell
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = $newStatus
$wsZh.Range("K3").Value = $zhDetail

# de-de detail sheet: row 3 status + new Error Detail cell
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = $newStatus
$wsDe.Range("K3").Value = $deDetail
